$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Neodymium")
$ws.Range("B2").Value = [double]"1.455684284667581E-17"
$ws.Range("C2").Value = [double]"0.003234355465845292"
$ws.Range("D2").Value = [double]"0.7019255690660977"
$ws.Range("E2").Value = [double]"0.9034683182445059"
$ws.Range("B3").Value = [double]"1.573017879546552E-07"
$ws.Range("C3").Value = [double]"0.0548952861088679"
$ws.Range("D3").Value = [double]"0.6513003445992317"
$ws.Range("E3").Value = [double]"0.8388692114087301"
$ws.Range("B4").Value = [double]"2.455525622174766E-09"
$ws.Range("C4").Value = [double]"0.05091945878934708"
$ws.Range("D4").Value = [double]"0.5322594535634154"
$ws.Range("E4").Value = [double]"0.6745337844781117"
$ws.Range("C5").Value = [double]"0.0001125284018505872"
$ws.Range("D5").Value = [double]"0.03432694848831427"
$ws.Range("E5").Value = [double]"0.05001255831081258"

$ws = $wb.Worksheets.Item("Dysprosium")
$ws.Range("B2").Value = [double]"1.455684284667678E-17"
$ws.Range("C2").Value = [double]"0.00323435546584551"
$ws.Range("D2").Value = [double]"0.7019255690661447"
$ws.Range("E2").Value = [double]"0.9034683182445664"
$ws.Range("B3").Value = [double]"1.573017879546658E-07"
$ws.Range("C3").Value = [double]"0.05489528610887159"
$ws.Range("D3").Value = [double]"0.6513003445992755"
$ws.Range("E3").Value = [double]"0.8388692114087867"
$ws.Range("B4").Value = [double]"2.455525622174932E-09"
$ws.Range("C4").Value = [double]"0.05091945878935049"
$ws.Range("D4").Value = [double]"0.5322594535634512"
$ws.Range("E4").Value = [double]"0.674533784478157"
$ws.Range("C5").Value = [double]"0.0001125284018505962"
$ws.Range("D5").Value = [double]"0.03432694848831701"
$ws.Range("E5").Value = [double]"0.05001255831081657"

$ws = $wb.Worksheets.Item("Copper")
$ws.Range("B2").Value = [double]"0.0009992687976177611"
$ws.Range("C2").Value = [double]"0.08229410131195405"
$ws.Range("D2").Value = [double]"1.58416514357327"
$ws.Range("E2").Value = [double]"2.070245676915336"
$ws.Range("B3").Value = [double]"0.007378544841876201"
$ws.Range("C3").Value = [double]"0.09549128864548032"
$ws.Range("D3").Value = [double]"1.113102957809399"
$ws.Range("E3").Value = [double]"1.552184907137255"
$ws.Range("B4").Value = [double]"0.01863605446332176"
$ws.Range("C4").Value = [double]"0.07341341649216067"
$ws.Range("D4").Value = [double]"1.325054156185886"
$ws.Range("E4").Value = [double]"1.788123587539221"
$ws.Range("B5").Value = [double]"0.006033682714979844"
$ws.Range("C5").Value = [double]"0.06910107071849605"
$ws.Range("D5").Value = [double]"1.248721570152209"
$ws.Range("E5").Value = [double]"1.757829693982964"

$ws = $wb.Worksheets.Item("Raw silicon")
$ws.Range("B2").Value = [double]"0.008341175675907163"
$ws.Range("C2").Value = [double]"0.1101398272791471"
$ws.Range("D2").Value = [double]"2.423778873971562"
$ws.Range("E2").Value = [double]"2.836977175650993"
$ws.Range("B3").Value = [double]"0.00966641920525574"
$ws.Range("C3").Value = [double]"0.1041167508083667"
$ws.Range("D3").Value = [double]"1.323147715383524"
$ws.Range("E3").Value = [double]"1.706756239850931"
$ws.Range("B4").Value = [double]"0.0527342217923131"
$ws.Range("C4").Value = [double]"0.09318616779647268"
$ws.Range("D4").Value = [double]"1.78872485157633"
$ws.Range("E4").Value = [double]"2.365167495720378"
$ws.Range("B5").Value = [double]"0.02918342513006148"
$ws.Range("C5").Value = [double]"0.107866079942882"
$ws.Range("D5").Value = [double]"2.136331103838482"
$ws.Range("E5").Value = [double]"2.530940646039109"

Write-Output "done"
